$d = $word.ActiveDocument

# Locate the paragraph range to replace: from the start of the "Завдання" heading
# paragraph through the end of the "Текст висновків" paragraph.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startPara -eq $null -and $t.StartsWith("Завдання")) {
        $startPara = $i
    }
    if ($t.StartsWith("Текст висновків")) {
        $endPara = $i
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate target paragraphs (start=$startPara end=$endPara)"
}

$rangeStart = $d.Paragraphs.Item($startPara).Range.Start
$rangeEnd = $d.Paragraphs.Item($endPara).Range.End
$target = $d.Range($rangeStart, $rangeEnd)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w14:paraId="205BFE66" w14:textId="628C3157" w:rsidR="00F47FDB" w:rsidRPr="00693266" w:rsidRDefault="00124A1F" w:rsidP="001F3D82"><w:pPr><w:pStyle w:val="Quote"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i w:val="0"/><w:iCs w:val="0"/></w:rPr><w:lastRenderedPageBreak/><w:t>Так бо БОГ</w:t></w:r><w:r w:rsidR="00693266" w:rsidRPr="00492233"><w:rPr><w:b/><w:bCs/><w:i w:val="0"/><w:iCs w:val="0"/></w:rPr><w:t xml:space="preserve"> полюбив світ, що дав Сина Свого Єдиного, аби кожен, хто вірує в Нього, не згинув, але мав життя вічне</w:t></w:r><w:r w:rsidR="00492233"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00693266" w:rsidRPr="008B213B"><w:rPr><w:i w:val="0"/><w:iCs w:val="0"/></w:rPr><w:t>(</w:t></w:r><w:hyperlink r:id="rId6" w:history="1"><w:r w:rsidR="00693266" w:rsidRPr="008B213B"><w:rPr><w:rStyle w:val="Hyperlink"/><w:i w:val="0"/><w:iCs w:val="0"/></w:rPr><w:t>Йоан 3:16</w:t></w:r></w:hyperlink><w:r w:rsidR="00693266" w:rsidRPr="008B213B"><w:rPr><w:i w:val="0"/><w:iCs w:val="0"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="P"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="679E726D" w14:textId="326BDC32" w:rsidR="00BE1B5B" w:rsidRDefault="00BE1B5B" w:rsidP="00BE1B5B"><w:pPr><w:pStyle w:val="H2"/></w:pPr><w:r><w:t>Завдання</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="14917D99" w14:textId="6EA826F5" w:rsidR="008F6190" w:rsidRDefault="008F6190" w:rsidP="00A13268"><w:pPr><w:pStyle w:val="H2"/></w:pPr><w:r><w:t>Висновки</w:t></w:r></w:p>
'@

$target.InsertXML($newXml) | Out-Null

Write-Output "done"
